$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits right after "syntax"
#    (it will be re-created further down in step 3, matching the
#    diff's relocation of the _GoBack marker).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) The "_Hlk20685137" bookmark keeps its name; its numeric w:id
#    attribute is renumbered automatically by the engine once the
#    lower-numbered "_GoBack" bookmark above is removed.
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# 3) Split "...hits large the blocks..." into "...hits large" + ","
#    + " the blocks..." (inserting a comma) and drop a fresh,
#    collapsed "_GoBack" bookmark right after the comma.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("hits large the blocks", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target sentence for edit"
}

# Position right after "large"
$afterLarge = $rng.Start + ("hits large").Length
$ip = $d.Range($afterLarge, $afterLarge)
$ip.InsertAfter(",")

# Force the newly inserted "," onto its own run (distinct <w:r>) by
# toggling a character attribute on/off - this splits the run boundary
# without leaving any lasting formatting difference behind.
$commaRng = $d.Range($afterLarge, $afterLarge + 1)
$commaRng.Bold = $true
$commaRng.Bold = $false

# Drop a fresh "_GoBack" bookmark (collapsed) right after the comma -
# Bookmarks.Add replaces any existing bookmark of the same name, so
# this also covers the removal from step 1 if it hadn't run.
$bmPoint = $d.Range($afterLarge + 1, $afterLarge + 1)
$d.Bookmarks.Add("_GoBack", $bmPoint) | Out-Null
